$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that flip from 0 to 1 per the recorded diff.
$cellsToSet = @(
    "G3", "H3",
    "G4", "H4",
    "D5", "E5",
    "D6", "E6",
    "D7", "E7",
    "H8",
    "D9", "E9",
    "H10",
    "D11", "E11",
    "H12",
    "H13",
    "H14",
    "H15",
    "H16",
    "H17",
    "H18"
)

foreach ($addr in $cellsToSet) {
    $ws.Range($addr).Value = 1
}
